$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C1").Value = "Q3"
$ws.Range("D1").Value = "Q4"

$ws.Range("D1").Select()
